$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" column header in H1, matching the style used by the other headers (s="1")
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in Save column values for rows 2-11
$saveValues = @(0, 1, 0, 0, 0, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
